# Flores: Removida bodega (se compran según necesidad)
# Delete the "Bodega" column (G), shift remaining columns left, rename
# "Cantidad Stock" header to "Cantidad Stock Actual", and fix up column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flores")

# Remove the entire "Bodega" column (G). This shifts H->G (Unidad) and
# I->H (Ultima Actualizacion) left by one, carrying their content/styles.
$ws.Columns("G").Delete()

# Rename the "Cantidad Stock" header (column F) to "Cantidad Stock Actual".
$ws.Range("F1").Value = "Cantidad Stock Actual"

# Adjust the column widths to match the regenerated layout. (ColumnWidth is
# expressed in characters and the engine adds ~5/6 of a character of
# padding when serializing to the OOXML "width" attribute, so we subtract
# it here to land on the clean target widths of 18 / 10 / 18.)
$pad = 5 / 6
$ws.Columns("F").ColumnWidth = 18 - $pad
$ws.Columns("G").ColumnWidth = 10 - $pad
$ws.Columns("H").ColumnWidth = 18 - $pad
